$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Scanner" to "Morgue"
$ws.Name = "Morgue"

# Update header E1: "Type" -> "Number"
$ws.Range("E1").Value = "Number"

# Ensure the data cells remain plain text (matches t="str" in the XML),
# not auto-converted to Excel dates/numbers.
$ws.Range("C2:E2").NumberFormat = "@"

$ws.Range("C2").Value = "04/05/2025"
$ws.Range("D2").Value = "15:48:38"
$ws.Range("E2").Value = "1746362918895"
